# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect newly generated output data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): rows 4, 6, 7
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3689
$ws1.Range("F6").Value = 37
$ws1.Range("F7").Value = 188

# Sheet "全部类型" (sheet4): rows 8, 10, 12
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F8").Value = 3689
$ws4.Range("F10").Value = 37
$ws4.Range("F12").Value = 188
